$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Snapshot the current Item/Half/Full/Image rows (rows 5-18) before moving anything.
# Index 0 in $data corresponds to worksheet row 5, index 13 corresponds to row 18.
$data = @()
for ($r = 5; $r -le 18; $r++) {
    $row = @($ws.Cells.Item($r, 1).Value2, $ws.Cells.Item($r, 3).Value2, $ws.Cells.Item($r, 4).Value2)
    $data += , $row
}

# The last three rows (old rows 16-18, the "Medium" pizza items) move up to
# immediately follow row 4, pushing the remaining rows (old rows 5-15) down.
$order = @(11, 12, 13, 0, 1, 2, 3, 4, 5, 6, 7, 8, 9, 10)

for ($i = 0; $i -lt $order.Length; $i++) {
    $newRow = 5 + $i
    $src = $data[$order[$i]]
    $ws.Cells.Item($newRow, 1).Value2 = $src[0]
    $ws.Cells.Item($newRow, 3).Value2 = $src[1]
    $ws.Cells.Item($newRow, 4).Value2 = $src[2]
}

# Match the author's final selection.
$ws.Range("A14").Select() | Out-Null
